$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap D2/D4 (dates) and M2/M4 (volumes) between row 2 and row 4
$ws.Range("D2").Value = 44874
$ws.Range("D4").Value = 44875

$ws.Range("M2").Value = 67
$ws.Range("M4").Value = 50
